# Regenerate s_vals data to filter save games.
# Updates the B:E (and derived G = sum(B:E)) columns on rows 2-6
# with the newly-computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    3 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    4 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    5 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    6 = @(0.04172184405617529, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($row, 2).Value = $b   # B
    $ws.Cells.Item($row, 3).Value = $c   # C
    $ws.Cells.Item($row, 4).Value = $d   # D
    $ws.Cells.Item($row, 5).Value = $e   # E
    $ws.Cells.Item($row, 7).Value = ($b + $c + $d + $e)   # G = sum(B:E)
}
